$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the BP1/BQ1 header labels (average_doctor <-> average_doctor_old)
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Update recomputed statistic values (rows 4-13)
$ws.Range("E4").Value = 0.42
$ws.Range("F4").Value = 0.071
$ws.Range("G4").Value = 0.266
$ws.Range("N4").Value = 0.422
$ws.Range("O4").Value = 0.06
$ws.Range("P4").Value = 0.244
$ws.Range("Q4").Value = 0.024
$ws.Range("R4").Value = 0.017
$ws.Range("S4").Value = 0.131
$ws.Range("W4").Value = 0.295
$ws.Range("X4").Value = 0.11
$ws.Range("Y4").Value = 0.331
$ws.Range("AI4").Value = 0.301
$ws.Range("AJ4").Value = 0.088
$ws.Range("AK4").Value = 0.297
$ws.Range("AU4").Value = 0.19
$ws.Range("AV4").Value = 0.029
$ws.Range("AW4").Value = 0.17
$ws.Range("BA4").Value = 1.994
$ws.Range("BB4").Value = 0.158
$ws.Range("BC4").Value = 0.397
$ws.Range("BG4").Value = 0.731
$ws.Range("BH4").Value = 0.139
$ws.Range("BI4").Value = 0.372
$ws.Range("BM4").Value = 0.714
$ws.Range("BN4").Value = 0.079
$ws.Range("BO4").Value = 0.281
$ws.Range("BP4").Value = 0.665
$ws.Range("BQ4").Value = 0.703
$ws.Range("E5").Value = 0.543
$ws.Range("F5").Value = 0.088
$ws.Range("G5").Value = 0.297
$ws.Range("N5").Value = 0.746
$ws.Range("O5").Value = 0.08
$ws.Range("P5").Value = 0.283
$ws.Range("Q5").Value = 0.016
$ws.Range("R5").Value = 0.007
$ws.Range("S5").Value = 0.084
$ws.Range("W5").Value = 0.285
$ws.Range("X5").Value = 0.11
$ws.Range("Y5").Value = 0.332
$ws.Range("AI5").Value = 0.323
$ws.Range("AJ5").Value = 0.098
$ws.Range("AK5").Value = 0.314
$ws.Range("AU5").Value = 0.368
$ws.Range("AV5").Value = 0.096
$ws.Range("AW5").Value = 0.309
$ws.Range("BA5").Value = 1.344
$ws.Range("BB5").Value = 0.082
$ws.Range("BC5").Value = 0.286
$ws.Range("BG5").Value = 0.402
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.225
$ws.Range("BM5").Value = 0.554
$ws.Range("BN5").Value = 0.065
$ws.Range("BO5").Value = 0.255
$ws.Range("BP5").Value = 0.448
$ws.Range("BQ5").Value = 0.454
$ws.Range("E6").Value = 0.474
$ws.Range("N6").Value = 0.539
$ws.Range("Q6").Value = 0.019
$ws.Range("W6").Value = 0.29
$ws.Range("AI6").Value = 0.312
$ws.Range("AU6").Value = 0.251
$ws.Range("BA6").Value = 1.598
$ws.Range("BG6").Value = 0.519
$ws.Range("BM6").Value = 0.624
$ws.Range("BP6").Value = 0.533
$ws.Range("BQ6").Value = 0.549
$ws.Range("E7").Value = 0.513
$ws.Range("N7").Value = 0.647
$ws.Range("Q7").Value = 0.017
$ws.Range("W7").Value = 0.287
$ws.Range("AI7").Value = 0.318
$ws.Range("AU7").Value = 0.31
$ws.Range("BA7").Value = 1.434
$ws.Range("BG7").Value = 0.442
$ws.Range("BM7").Value = 0.58
$ws.Range("BP7").Value = 0.478
$ws.Range("BQ7").Value = 0.488
$ws.Range("E8").Value = 0.604
$ws.Range("F8").Value = 0.113
$ws.Range("G8").Value = 0.336
$ws.Range("N8").Value = 0.775
$ws.Range("O8").Value = 0.068
$ws.Range("P8").Value = 0.26
$ws.Range("Q8").Value = 0.018
$ws.Range("W8").Value = 0.314
$ws.Range("X8").Value = 0.121
$ws.Range("Y8").Value = 0.348
$ws.Range("AI8").Value = 0.345
$ws.Range("AJ8").Value = 0.129
$ws.Range("AK8").Value = 0.36
$ws.Range("AU8").Value = 0.309
$ws.Range("AV8").Value = 0.087
$ws.Range("AW8").Value = 0.295
$ws.Range("BA8").Value = 1.743
$ws.Range("BB8").Value = 0.126
$ws.Range("BC8").Value = 0.355
$ws.Range("BG8").Value = 0.567
$ws.Range("BH8").Value = 0.106
$ws.Range("BI8").Value = 0.325
$ws.Range("BM8").Value = 0.692
$ws.Range("BN8").Value = 0.067
$ws.Range("BO8").Value = 0.259
$ws.Range("BP8").Value = 0.581
$ws.Range("BQ8").Value = 0.602
$ws.Range("E9").Value = 0.544
$ws.Range("F9").Value = 0.248
$ws.Range("G9").Value = 0.498
$ws.Range("N9").Value = 0.678
$ws.Range("O9").Value = 0.218
$ws.Range("P9").Value = 0.467
$ws.Range("W9").Value = 0.211
$ws.Range("X9").Value = 0.167
$ws.Range("Y9").Value = 0.408
$ws.Range("AI9").Value = 0.267
$ws.Range("AJ9").Value = 0.196
$ws.Range("AK9").Value = 0.442
$ws.Range("BA9").Value = 1.688
$ws.Range("BB9").Value = 0.247
$ws.Range("BC9").Value = 0.497
$ws.Range("BG9").Value = 0.6
$ws.Range("BH9").Value = 0.24
$ws.Range("BI9").Value = 0.49
$ws.Range("BM9").Value = 0.644
$ws.Range("BN9").Value = 0.229
$ws.Range("BO9").Value = 0.479
$ws.Range("BP9").Value = 0.563
$ws.Range("BQ9").Value = 0.583
$ws.Range("E10").Value = 0.678
$ws.Range("F10").Value = 0.218
$ws.Range("G10").Value = 0.467
$ws.Range("N10").Value = 0.867
$ws.Range("O10").Value = 0.116
$ws.Range("P10").Value = 0.34
$ws.Range("W10").Value = 0.389
$ws.Range("X10").Value = 0.238
$ws.Range("Y10").Value = 0.487
$ws.Range("AI10").Value = 0.378
$ws.Range("AJ10").Value = 0.235
$ws.Range("AK10").Value = 0.485
$ws.Range("AU10").Value = 0.3
$ws.Range("AV10").Value = 0.21
$ws.Range("AW10").Value = 0.458
$ws.Range("BA10").Value = 2.078
$ws.Range("BB10").Value = 0.244
$ws.Range("BC10").Value = 0.494
$ws.Range("BG10").Value = 0.656
$ws.Range("BH10").Value = 0.226
$ws.Range("BI10").Value = 0.475
$ws.Range("BM10").Value = 0.844
$ws.Range("BN10").Value = 0.131
$ws.Range("BO10").Value = 0.362
$ws.Range("BP10").Value = 0.693
$ws.Range("BQ10").Value = 0.725
$ws.Range("E11").Value = 0.711
$ws.Range("F11").Value = 0.205
$ws.Range("G11").Value = 0.453
$ws.Range("N11").Value = 0.889
$ws.Range("O11").Value = 0.099
$ws.Range("P11").Value = 0.314
$ws.Range("W11").Value = 0.389
$ws.Range("X11").Value = 0.238
$ws.Range("Y11").Value = 0.487
$ws.Range("AI11").Value = 0.411
$ws.Range("AJ11").Value = 0.242
$ws.Range("AK11").Value = 0.492
$ws.Range("AU11").Value = 0.422
$ws.Range("AV11").Value = 0.244
$ws.Range("AW11").Value = 0.494
$ws.Range("BA11").Value = 2.078
$ws.Range("BB11").Value = 0.244
$ws.Range("BC11").Value = 0.494
$ws.Range("BG11").Value = 0.656
$ws.Range("BH11").Value = 0.226
$ws.Range("BI11").Value = 0.475
$ws.Range("BM11").Value = 0.844
$ws.Range("BN11").Value = 0.131
$ws.Range("BO11").Value = 0.362
$ws.Range("BP11").Value = 0.693
$ws.Range("BQ11").Value = 0.727
$ws.Range("E12").Value = 1.422
$ws.Range("F12").Value = 0.775
$ws.Range("G12").Value = 0.88
$ws.Range("N12").Value = 1.476
$ws.Range("O12").Value = 1.079
$ws.Range("P12").Value = 1.039
$ws.Range("W12").Value = 1.629
$ws.Range("X12").Value = 0.576
$ws.Range("Y12").Value = 0.759
$ws.Range("AI12").Value = 1.703
$ws.Range("AJ12").Value = 1.29
$ws.Range("AK12").Value = 1.136
$ws.Range("AU12").Value = 2.7
$ws.Range("AV12").Value = 2.76
$ws.Range("AW12").Value = 1.661
$ws.Range("BA12").Value = 3.732
$ws.Range("BB12").Value = 0.412
$ws.Range("BC12").Value = 0.642
$ws.Range("BG12").Value = 1.102
$ws.Range("BH12").Value = 0.125
$ws.Range("BI12").Value = 0.354
$ws.Range("BM12").Value = 1.303
$ws.Range("BN12").Value = 0.343
$ws.Range("BO12").Value = 0.585
$ws.Range("BP12").Value = 1.244
$ws.Range("BQ12").Value = 1.267
$ws.Range("E13").Value = 1.595
$ws.Range("F13").Value = 0.66
$ws.Range("G13").Value = 0.812
$ws.Range("N13").Value = 2.113
$ws.Range("O13").Value = 0.9
$ws.Range("P13").Value = 0.948
$ws.Range("W13").Value = 1.049
$ws.Range("X13").Value = 0.196
$ws.Range("Y13").Value = 0.442
$ws.Range("AI13").Value = 1.277
$ws.Range("AJ13").Value = 0.374
$ws.Range("AK13").Value = 0.612
$ws.Range("AU13").Value = 2.284
$ws.Range("AV13").Value = 0.949
$ws.Range("AW13").Value = 0.974
$ws.Range("BA13").Value = 2.367
$ws.Range("BB13").Value = 0.297
$ws.Range("BC13").Value = 0.545
$ws.Range("BG13").Value = 0.59
$ws.Range("BH13").Value = 0.072
$ws.Range("BI13").Value = 0.269
$ws.Range("BM13").Value = 0.901
$ws.Range("BN13").Value = 0.282
$ws.Range("BO13").Value = 0.531
$ws.Range("BP13").Value = 0.789
$ws.Range("BQ13").Value = 0.729
